{"js": "// Rewrote IB work experience bullets (and a couple of Education bullets)\n// to match the new resume copy. Each edit is a straight text replacement\n// of an existing paragraph's content \u2014 the paragraph count/order is\n// unchanged, only the wording changes \u2014 so we find each paragraph by its\n// distinctive (stable) original leading text and replace its full text.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nfunction findIndex(items, startsWith) {\n  for (let i = 0; i < items.length; i++) {\n    if (items[i].text.indexOf(startsWith) === 0) {\n      return i;\n    }\n  }\n  throw new Error(\"Paragraph not found: \" + startsWith);\n}\n\nconst items = paragraphs.items;\n\nconst replacements = [\n  [\n    \"Significantly increase regulatory reporting efficiency\",\n    \"Design and implement six regulatory reporting systems handling millions of records daily in less than a year in Python and Go often with three weeks or less until go-live with a 100% on time record\",\n  ],\n  [\n    \"Improve code readability and maintainability\",\n    \"Introduce the latest technologies and best-practices to legacy systems by redesigning projects written in Perl and implementing them using object-oriented design patterns in Python\",\n  ],\n  [\n    \"Facilitate relationships between Interactive Brokers\",\n    \"Work directly with compliance analysts to automate and simplify delivery of customer trading data to affiliated broker dealers, allowing for significantly more efficient bookkeeping for both Interactive Brokers and the affiliated brokers\",\n  ],\n  [\n    \"Developed a full-scale e-commerce site\",\n    \"Independently developed a full-scale e-commerce site using PHP and MySQL by implementing software development skills, including database design and object-oriented design, acquired through professional experience and theoretical course work\",\n  ],\n  [\n    \"Developed web applications using web technologies\",\n    \"Created professional-grade web applications using JavaScript, jQuery, Bootstrap, and Flask\",\n  ],\n];\n\nfor (const [needle, newText] of replacements) {\n  const idx = findIndex(items, needle);\n  items[idx].insertText(newText, \"Replace\");\n}\n\nawait context.sync();\n", "ps1": "# Rewrote IB work experience bullets (and a couple of Education bullets)\n# to match the new resume copy. Each edit is a straight text replacement\n# of an existing paragraph's content (paragraph count/order is unchanged,\n# only the wording changes), so we drive it with Find/Replace over the\n# whole document, one exact phrase at a time.\n\n$d = $word.ActiveDocument\n\nfunction Replace-ExactText($oldText, $newText) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $newText\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n\nReplace-ExactText `\n    \"Significantly increase regulatory reporting efficiency by designing and implementing systems in Python and Go to process data from Oracle databases as well as CSV and XML files for delivery to US and foreign government agencies\" `\n    \"Design and implement six regulatory reporting systems handling millions of records daily in less than a year in Python and Go often with three weeks or less until go-live with a 100% on time record\"\n\nReplace-ExactText `\n    \"Improve code readability and maintainability by redesigning systems written in Perl and implementing them using object-oriented design patterns in Python\" `\n    \"Introduce the latest technologies and best-practices to legacy systems by redesigning projects written in Perl and implementing them using object-oriented design patterns in Python\"\n\nReplace-ExactText `\n    \"Facilitate relationships between Interactive Brokers and affiliated broker-dealers by automating and simplifying delivery of customer trading data, allowing for more efficient bookkeeping on both sides\" `\n    \"Work directly with compliance analysts to automate and simplify delivery of customer trading data to affiliated broker dealers, allowing for significantly more efficient bookkeeping for both Interactive Brokers and the affiliated brokers\"\n\nReplace-ExactText `\n    \"Developed a full-scale e-commerce site using PHP and MySQL by implementing software development skills, including database design and object-oriented design, acquired through professional experience and theoretical course work\" `\n    \"Independently developed a full-scale e-commerce site using PHP and MySQL by implementing software development skills, including database design and object-oriented design, acquired through professional experience and theoretical course work\"\n\nReplace-ExactText `\n    \"Developed web applications using web technologies including JavaScript, jQuery, Bootstrap, and Flask\" `\n    \"Created professional-grade web applications using JavaScript, jQuery, Bootstrap, and Flask\"\n"}
